$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Price column (D) updates that look numeric -> force text to match source formatting
Set-TextValue "D5" "208.69"
Set-TextValue "D15" "0.534"
Set-TextValue "D17" "63.58"
Set-TextValue "D18" "219.73"
Set-TextValue "D22" "4.17"
Set-TextValue "D23" "9.77"
Set-TextValue "D24" "1.98"
Set-TextValue "D25" "154.44"
Set-TextValue "D26" "6.76"
Set-TextValue "D28" "15.18"
Set-TextValue "D31" "0.0473"
Set-TextValue "D35" "1.54"
Set-TextValue "D39" "0.538"
Set-TextValue "D40" "0.829"
Set-TextValue "D42" "0.978"
Set-TextValue "D43" "64.70"
Set-TextValue "D48" "87.03"

# Remaining text updates (prices with multiple dots, and all Volume(1h) percentage cells)
$ws.Range("D2").Value = "27.734.23"
$ws.Range("D3").Value = "1.596.26"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "1.823.54"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "1.588.70"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("E14").Value = "  -3.37%  "
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("D16").Value = "27.731.05"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("E18").Value = "  -3.49%  "
$ws.Range("D19").Value = "0.0₃0699"
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  -4.20%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  -4.52%  "
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("E32").Value = "  -4.57%  "
$ws.Range("D33").Value = "1.379.46"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("E35").Value = "  -4.04%  "
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("E45").Value = "  -2.65%  "
$ws.Range("E46").Value = "  -3.64%  "
$ws.Range("D47").Value = "1.733.82"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  -3.42%  "
$ws.Range("E51").Value = "  -1.19%  "
